$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H; this shifts the old H -> I and old I -> J,
# carrying over cell formatting/styles automatically.
$ws.Columns.Item(8).Insert()

# New column header (new shared string "prUd.5,e.0")
$ws.Cells.Item(1,8).Value = "prUd.5,e.0"

# The Insert() leaves the new H1 without the label style used by the
# other header cells (s="1", quotePrefix) - copy that formatting over from G1.
$ws.Cells.Item(1,7).Copy()
$ws.Cells.Item(1,8).PasteSpecial(-4122)

# Data rows 2-28: new/updated results for columns H (new "prUd.5,e.0" run),
# I ("Itax=0.1", recomputed) and J ("Otax=0.1", recomputed)
$ws.Cells.Item(2,8).Value = 1.0367159805528305
$ws.Cells.Item(2,9).Value = 1.0218924908754679
$ws.Cells.Item(2,10).Value = 0.98256190948137989
$ws.Cells.Item(3,8).Value = 1.0427805239888699
$ws.Cells.Item(3,9).Value = 1.0641978434498767
$ws.Cells.Item(3,10).Value = 1.1228202670038634
$ws.Cells.Item(4,8).Value = 1.03915507046988
$ws.Cells.Item(4,9).Value = 1.0387071207235727
$ws.Cells.Item(4,10).Value = 1.0365510187680878
$ws.Cells.Item(5,8).Value = 1.0023527079378749
$ws.Cells.Item(5,9).Value = 1.0164544020024273
$ws.Cells.Item(5,10).Value = 1.0549472850165791
$ws.Cells.Item(6,8).Value = 0.99573648804003345
$ws.Cells.Item(6,9).Value = 0.9707623117768589
$ws.Cells.Item(6,10).Value = 0.9073237033949143
$ws.Cells.Item(7,8).Value = 1
$ws.Cells.Item(7,9).Value = 1
$ws.Cells.Item(7,10).Value = 1
$ws.Cells.Item(8,8).Value = 0.94435493123749292
$ws.Cells.Item(8,9).Value = 0.89431602404483068
$ws.Cells.Item(8,10).Value = 0.83994765227154744
$ws.Cells.Item(9,8).Value = 1.0388474003327777
$ws.Cells.Item(9,9).Value = 1.0365981349629205
$ws.Cells.Item(9,10).Value = 0.96509548914632604
$ws.Cells.Item(10,8).Value = 80
$ws.Cells.Item(10,9).Value = 80
$ws.Cells.Item(10,10).Value = 80
$ws.Cells.Item(11,8).Value = 54
$ws.Cells.Item(11,9).Value = 54
$ws.Cells.Item(11,10).Value = 54
$ws.Cells.Item(12,8).Value = 124
$ws.Cells.Item(12,9).Value = 124
$ws.Cells.Item(12,10).Value = 124
$ws.Cells.Item(13,8).Value = 31.842456944360443
$ws.Cells.Item(13,9).Value = 30.997413605068637
$ws.Cells.Item(13,10).Value = 30.828302885539017
$ws.Cells.Item(14,8).Value = 48.243500807548365
$ws.Cells.Item(14,9).Value = 49.028373084915735
$ws.Cells.Item(14,10).Value = 49.189565550386192
$ws.Cells.Item(15,8).Value = 25.305819795566755
$ws.Cells.Item(15,9).Value = 26.051524132661086
$ws.Cells.Item(15,10).Value = 25.925150004988197
$ws.Cells.Item(16,8).Value = 28.755036236909774
$ws.Cells.Item(16,9).Value = 28.094657293907357
$ws.Cells.Item(16,10).Value = 28.204163637656286
$ws.Cells.Item(17,8).Value = 79.812225144363239
$ws.Cells.Item(17,9).Value = 78.704956997971451
$ws.Cells.Item(17,10).Value = 75.833173027923138
$ws.Cells.Item(18,8).Value = 44.188397762351457
$ws.Cells.Item(18,9).Value = 45.325204188719908
$ws.Cells.Item(18,10).Value = 48.494269283791553
$ws.Cells.Item(19,8).Value = 139.2024749421293
$ws.Cells.Item(19,9).Value = 138.88306022726826
$ws.Cells.Item(19,10).Value = 137.93719941830233
$ws.Cells.Item(20,8).Value = 128.85522873815168
$ws.Cells.Item(20,9).Value = 128.79968296955442
$ws.Cells.Item(20,10).Value = 128.53232632723834
$ws.Cells.Item(21,8).Value = 10.391550704689651
$ws.Cells.Item(21,9).Value = 10.387071207222132
$ws.Cells.Item(21,10).Value = 10.365510187680512
$ws.Cells.Item(22,8).Value = 1.0391550704689652
$ws.Cells.Item(22,9).Value = 1.0387071207222132
$ws.Cells.Item(22,10).Value = 1.0365510187680511
$ws.Cells.Item(23,8).Value = 1
$ws.Cells.Item(23,9).Value = 1
$ws.Cells.Item(23,10).Value = 1
$ws.Cells.Item(24,8).Value = 0.9933993095988608
$ws.Cells.Item(24,9).Value = 0.95504757504561499
$ws.Cells.Item(24,10).Value = 0.860065442398532
$ws.Cells.Item(25,8).Value = 0.99765281430454245
$ws.Cells.Item(25,9).Value = 0.9838119624746452
$ws.Cells.Item(25,10).Value = 0.94791466284903936
$ws.Cells.Item(26,8).Value = 0.94213835485145747
$ws.Cells.Item(26,9).Value = 0.87983880268806691
$ws.Cells.Item(26,10).Value = 0.79619869561382606
$ws.Cells.Item(27,8).Value = 1.0364090325749533
$ws.Cells.Item(27,9).Value = 1.019817645455428
$ws.Cells.Item(27,10).Value = 0.91482816521126842
$ws.Cells.Item(28,8).Value = 138.87574088417284
$ws.Cells.Item(28,9).Value = 136.63481603667313
$ws.Cells.Item(28,10).Value = 130.75269388094077
